$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 14 (shifts old rows 14..41 down to 15..42)
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly record
$ws.Range("A14").Value2 = 10
$ws.Range("B14").Value2 = "Vega Modelo de Temuco"
$ws.Range("C14").Value2 = "La Araucanía"
$ws.Range("D14").Value2 = 44477
$ws.Range("E14").Value2 = 9
$ws.Range("F14").Value2 = 300000001
$ws.Range("G14").Value2 = "Rabanito"
$ws.Range("H14").Value2 = "Sin especificar"
$ws.Range("I14").Value2 = "Primera"
$ws.Range("J14").Value2 = 20
$ws.Range("K14").Value2 = 8000
$ws.Range("L14").Value2 = 8000
$ws.Range("M14").Value2 = 8000
$ws.Range("N14").Value2 = '$/docena de paquetes'
$ws.Range("O14").Value2 = "Provincia de Cautín"
$ws.Range("P14").Value2 = 667
$ws.Range("Q14").Value2 = 12
$ws.Range("R14").Value2 = "Hortaliza"
